{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// --- Change 1: remove the \"Ready made drones or DIY kit drones...\" bullet paragraph entirely ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (\n    paragraphs.items[i].text.indexOf(\n      \"Ready made drones or DIY kit drones are not allowed to participate in the competition.\"\n    ) !== -1\n  ) {\n    paragraphs.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// --- Change 2: merge the split \"Penalties are imposed ...\" runs into a single run of text ---\nconst mergedText =\n  \"Penalties are imposed if the drone touches the border of the track or if skips or replacements are chosen. Each penalty adds extra time to the total time. \";\n\nconst searchResults = body.search(mergedText, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(mergedText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: remove the \"Ready made drones or DIY kit drones...\" bullet paragraph entirely ---\n$targetText = \"Ready made drones or DIY kit drones are not allowed to participate in the competition.\"\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$targetText*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- Change 2: merge the split \"Penalties are imposed ...\" runs into a single run of text ---\n# In the source document this sentence is split across three separate runs\n# (\"...skips or\" + \" \" + \"replacements...\"); Find/Replace operates on the\n# paragraph's logical text stream, so it matches across those run boundaries\n# and replaces them with one run containing the full sentence.\n$searchText = \"Penalties are imposed if the drone touches the border of the track or if skips or replacements are chosen. Each penalty adds extra time to the total time. \"\n$mergedText = $searchText\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    $searchText,   # FindText\n    $false,        # MatchCase\n    $false,        # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    1,             # Wrap            (wdFindContinue)\n    $false,        # Format\n    $mergedText,   # ReplaceWith\n    2              # Replace         (wdReplaceAll)\n)\n"}
